$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Split the run containing the job-dispatcher name so that
#    "allToLargest" sits in its own run (the quoted name is typed
#    as its own segment, matching the proofed word boundary).
# -----------------------------------------------------------------
$full = $d.Content.Text
$nameWord = "allToLargest"
$idx = $full.IndexOf($nameWord)
$nameRange = $d.Range($idx, $idx + $nameWord.Length)
$nameRange.Font.Bold = 1
$nameRange.Font.Bold = 0

# -----------------------------------------------------------------
# 2) Remove the extra blank paragraph that used to sit between the
#    Introduction paragraph and the "System Overview:" heading.
# -----------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim() -eq "" -and $p.Range.Font.Bold -eq 0 -and $p.Range.Font.Size -eq 14) {
        $p.Range.Delete()
        break
    }
}

# -----------------------------------------------------------------
# 3) Add a new paragraph right after "Implementation:" describing
#    the client classes.
# -----------------------------------------------------------------
# Build the new paragraph's text on a scratch paragraph split off of
# the Introduction body paragraph, since it already carries the
# plain (non-bold, 24-half-point) character formatting we need; a
# paragraph inserted next to the bold "Implementation:"/"References:"
# headings would otherwise inherit their bold run properties.
$introPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("The goal of this project")) {
        $introPara = $p
        break
    }
}
$introPara.Range.InsertParagraphAfter()
$scratchIndex = $introPara.Index + 1
$scratchPara = $d.Paragraphs.Item($scratchIndex)

$newText = "Our client is implemented using 3 classes. The ClientJobScheduler class, the Jobs class, and the Servers class. "
$scratchPara.Range.Text = $newText

# Split "ClientJobScheduler" into its own run within the scratch
# paragraph, matching the word boundary in the target document.
$scratchPara = $d.Paragraphs.Item($scratchIndex)
$pStart = $scratchPara.Range.Start
$classWord = "ClientJobScheduler"
$wordIdx = $newText.IndexOf($classWord)
$classRange = $d.Range($pStart + $wordIdx, $pStart + $wordIdx + $classWord.Length)
$classRange.Font.Bold = 1
$classRange.Font.Bold = 0

# Grab the scratch paragraph, mark included, as rich (formatted)
# text so its clean formatting travels with it.
$scratchPara = $d.Paragraphs.Item($scratchIndex)
$srcRange = $d.Range($scratchPara.Range.Start, $scratchPara.Range.End)
$richText = $srcRange.FormattedText

# Locate "Implementation:" and paste the rich text immediately
# after it.
$implPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Implementation:")) {
        $implPara = $p
        break
    }
}
$destPos = $implPara.Range.End
$destRange = $d.Range($destPos, $destPos)
$destRange.FormattedText = $richText

# Remove the scratch paragraph that was only used to build the
# cleanly-formatted rich text.
$scratchPara = $d.Paragraphs.Item($scratchIndex)
$scratchPara.Range.Delete()
